$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Rows 19 and 20 swap: Uniswap/Dai trade places (name, link, price, volume)
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D19" "0.9996"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "5.570"
$ws.Range("E20").Value = "  +6.35%  "

# Remaining price/volume updates
Set-TextValue "D2" "30.614.13"
$ws.Range("E2").Value = "  +0.74%  "
Set-TextValue "D3" "1.922.57"
$ws.Range("E3").Value = "  +0.08%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "246.78"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("E6").Value = "  +0.04%  "
Set-TextValue "D7" "0.4738"
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue "D8" "0.2884"
$ws.Range("E8").Value = "  +1.08%  "
Set-TextValue "D9" "0.06837"
$ws.Range("E9").Value = "  +3.70%  "
Set-TextValue "D10" "105.08"
$ws.Range("E10").Value = "  -0.65%  "
Set-TextValue "D11" "18.35"
$ws.Range("E11").Value = "  -4.12%  "
Set-TextValue "D12" "1.914.68"
$ws.Range("E12").Value = "  -0.33%  "
Set-TextValue "D13" "0.07693"
$ws.Range("E13").Value = "  +1.36%  "
Set-TextValue "D14" "5.331"
$ws.Range("E14").Value = "  +4.17%  "
Set-TextValue "D15" "0.6673"
$ws.Range("E15").Value = "  +1.74%  "
Set-TextValue "D16" "291.89"
$ws.Range("E16").Value = "  -4.21%  "
Set-TextValue "D17" "30.627.72"
$ws.Range("E17").Value = "  +0.70%  "
Set-TextValue "D18" "0.000007619"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E21").Value = "  +0.19%  "
Set-TextValue "D22" "2.174.50"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  -0.07%  "
Set-TextValue "D24" "6.447"
$ws.Range("E24").Value = "  +2.71%  "
Set-TextValue "D25" "9.452"
$ws.Range("E25").Value = "  +2.93%  "
Set-TextValue "D26" "167.60"
$ws.Range("E26").Value = "  +0.17%  "
Set-TextValue "D27" "21.02"
$ws.Range("E27").Value = "  +6.53%  "
Set-TextValue "D28" "2.115"
$ws.Range("E28").Value = "  +5.34%  "
Set-TextValue "D29" "0.1071"
$ws.Range("E29").Value = "  -4.63%  "
Set-TextValue "D30" "1.394"
$ws.Range("E30").Value = "  +3.71%  "
Set-TextValue "D31" "4.177"
$ws.Range("E31").Value = "  +1.87%  "
Set-TextValue "D32" "4.058"
$ws.Range("E32").Value = "  +3.62%  "
Set-TextValue "D33" "0.05041"
$ws.Range("E33").Value = "  +0.60%  "
Set-TextValue "D34" "0.7375"
$ws.Range("E34").Value = "  -0.33%  "
Set-TextValue "D35" "1.144"
$ws.Range("E35").Value = "  +0.06%  "
Set-TextValue "D36" "0.02066"
$ws.Range("E36").Value = "  +6.21%  "
Set-TextValue "D37" "2.741"
$ws.Range("E37").Value = "  +1.02%  "
Set-TextValue "D38" "2.690"
$ws.Range("E38").Value = "  -0.15%  "
Set-TextValue "D39" "2.055"
$ws.Range("E39").Value = "  +0.58%  "
Set-TextValue "D40" "111.27"
$ws.Range("E40").Value = "  +3.69%  "
Set-TextValue "D41" "0.8765"
$ws.Range("E41").Value = "  +0.30%  "
Set-TextValue "D42" "0.4373"
$ws.Range("E42").Value = "  +6.18%  "
Set-TextValue "D43" "5.913"
$ws.Range("E43").Value = "  +2.05%  "
Set-TextValue "D44" "1.0000"
$ws.Range("E44").Value = "  +0.09%  "
Set-TextValue "D45" "67.96"
$ws.Range("E45").Value = "  -2.00%  "
Set-TextValue "D46" "7.275"
$ws.Range("E46").Value = "  +1.05%  "
Set-TextValue "D47" "9.384"
$ws.Range("E47").Value = "  +1.54%  "
Set-TextValue "D48" "48.19"
$ws.Range("E48").Value = "  +14.37%  "
Set-TextValue "D49" "0.1242"
$ws.Range("E49").Value = "  +3.31%  "
Set-TextValue "D50" "0.2533"
$ws.Range("E50").Value = "  +13.46%  "
Set-TextValue "D51" "34.98"
$ws.Range("E51").Value = "  +0.70%  "
